$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.583.28"
$ws.Range("E2").Value = "  -0.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.875.13"
$ws.Range("E3").Value = "  -1.11%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.53"
$ws.Range("E5").Value = "  +0.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("E7").Value = "  -1.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2896"
$ws.Range("E8").Value = "  -0.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06472"
$ws.Range("E9").Value = "  -1.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.99"
$ws.Range("E10").Value = "  +2.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07725"
$ws.Range("E11").Value = "  -0.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7419"
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.872.88"
$ws.Range("E13").Value = "  -1.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "96.00"
$ws.Range("E14").Value = "  -0.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.167"
$ws.Range("E15").Value = "  -0.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "274.24"
$ws.Range("E16").Value = "  -3.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.651.02"
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.26"
$ws.Range("E18").Value = "  -3.32%  "
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007464"
$ws.Range("E20").Value = "  -2.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.118.34"
$ws.Range("E21").Value = "  -1.36%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.192"
$ws.Range("E23").Value = "  -2.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.167"
$ws.Range("E24").Value = "  -1.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.196"
$ws.Range("E25").Value = "  -1.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.00"
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("E27").Value = "  -2.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.902"
$ws.Range("E28").Value = "  -5.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.09933"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.347"
$ws.Range("E30").Value = "  -2.94%  "
$ws.Range("E31").Value = "  -0.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.234"
$ws.Range("E32").Value = "  -3.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.082"
$ws.Range("E33").Value = "  -1.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04760"
$ws.Range("E34").Value = "  -0.50%  "
$ws.Range("E35").Value = "  -1.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6912"
$ws.Range("E36").Value = "  -2.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.719"
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01843"
$ws.Range("E38").Value = "  -1.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.759"
$ws.Range("E39").Value = "  -0.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.254"
$ws.Range("E40").Value = "  -4.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.17"
$ws.Range("E41").Value = "  +2.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.958"
$ws.Range("E42").Value = "  +1.03%  "
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4154"
$ws.Range("E44").Value = "  -1.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8331"
$ws.Range("E45").Value = "  -2.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.27"
$ws.Range("E46").Value = "  -1.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.401"
$ws.Range("E47").Value = "  -0.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.32"
$ws.Range("E48").Value = "  -0.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.963"
$ws.Range("E49").Value = "  -2.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "912.02"
$ws.Range("E50").Value = "  -2.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05667"
$ws.Range("E51").Value = "  +0.83%  "
